# Add test-data rows for the "add deals" page to the "deals" worksheet (sheet3).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("deals")

# --- Header row (row 1) -----------------------------------------------
$headers = @("title","company","primarycontact","amount","probability","commission","identifier","tags","descrption","nextStep","product","quantity","type","source")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- Data row 2 (Flipkart deal) ----------------------------------------
$row2 = @("deals title - 1","Flipkart company","bansal - contact",5000,80,20,"test identifier","tagOne, tagTwo, tagThree","test desc -- added by salesperson","waiting for answer from client","Test Product",3,"Priority","Word of Mouth")
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}

# --- Data row 3 (Amazon deal) -------------------------------------------
$row3 = @("deals title - 2","Amazon Ccompany","John - Contact",8000,60,10,"Test - 2","tagFour, tagFive","test desc -- added by salesperson (amazon)","amazon - next step","Test Product",5,"Priority","Word of Mouth")
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, $i + 1).Value = $row3[$i]
}

# --- Header row styling (yellow fill, same as other sheets' headers) ----
$ws.Range("A1:N1").Interior.Color = 65535

# --- Column widths (best-fit, approximated to this engine's 1/6 step) ---
$ws.Columns.Item(1).ColumnWidth = 11.666666666666666
$ws.Columns.Item(2).ColumnWidth = 17.166666666666668
$ws.Columns.Item(3).ColumnWidth = 14
$ws.Columns.Item(5).ColumnWidth = 9.833333333333334
$ws.Columns.Item(6).ColumnWidth = 10.666666666666666
$ws.Columns.Item(7).ColumnWidth = 12.666666666666666
$ws.Columns.Item(8).ColumnWidth = 23.333333333333332
$ws.Columns.Item(9).ColumnWidth = 39.833333333333336
$ws.Columns.Item(10).ColumnWidth = 27.333333333333332
$ws.Columns.Item(11).ColumnWidth = 11.166666666666666
$ws.Columns.Item(12).ColumnWidth = 7.666666666666667
$ws.Columns.Item(13).ColumnWidth = 6.666666666666667
$ws.Columns.Item(14).ColumnWidth = 13.666666666666666

# --- Activate the deals sheet/tab and set its selection ------------------
$ws.Activate() | Out-Null
$ws.Range("M11").Select() | Out-Null

# --- Remove the tab selection / selection changes from "contacts" sheet --
# (Activating "deals" automatically clears tabSelected on the previously
# selected sheet; the "contacts" sheet's own selection (J7) is unchanged.)

Write-Host "deals test data added"
